$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the first three fixture rows (West Ham Carabao Cup, Aston Villa, Villarreal CF).
# Remaining rows shift up, turning the old 21-row list into 18 rows.
$ws.Range("A1:B3").EntireRow.Delete()

# The Manchester City fixture kick-off time moved from 15:00 to 12:30.
# After the deletion above, that fixture is now on row 4.
$ws.Range("B4").Value = "06 NovSat12:30"
